$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated correlation matrix values (recomputed macro section).
$cellUpdates = @(
    @("D2", [double]"0.04250259628683546"),
    @("E2", [double]"0.06052232935307845"),
    @("F2", [double]"-0.03737587636125757"),
    @("G2", [double]"0.003932712359607337"),
    @("H2", [double]"0.03251026699987748"),
    @("I2", [double]"0.1455275192583088"),
    @("J2", [double]"0.03110426779743944"),
    @("K2", [double]"-0.06786464654325448"),
    @("L2", [double]"0.09694398396020938"),
    @("M2", [double]"0.1693593265268756"),
    @("D3", [double]"0.02090081364037645"),
    @("E3", [double]"0.05405000766046579"),
    @("F3", [double]"0.1691006458289873"),
    @("G3", [double]"0.1131570554671354"),
    @("H3", [double]"0.03464555403083762"),
    @("I3", [double]"-0.05464967430273746"),
    @("J3", [double]"-0.02353409310941642"),
    @("K3", [double]"0.1086358250565714"),
    @("L3", [double]"0.1136429820580396"),
    @("M3", [double]"-0.03672462593629056"),
    @("B4", [double]"0.04250259628683546"),
    @("C4", [double]"0.02090081364037645"),
    @("E4", [double]"0.5636986541578768"),
    @("F4", [double]"0.1553023117281974"),
    @("G4", [double]"0.6735708521826163"),
    @("H4", [double]"0.9869596445666107"),
    @("I4", [double]"0.7703047989162173"),
    @("J4", [double]"0.6657197040299973"),
    @("K4", [double]"0.2524737509893859"),
    @("L4", [double]"0.4706509334439072"),
    @("M4", [double]"0.6081882455058323"),
    @("N4", [double]"3.566671208812168e-16"),
    @("O4", [double]"-0.1375382841619496"),
    @("Q4", [double]"4.048653804597595e-16"),
    @("R4", [double]"-0.1123992994262411"),
    @("S4", [double]"0.07052404406482074"),
    @("T4", [double]"0.2153517673687143"),
    @("U4", [double]"-0.1390079924904912"),
    @("B5", [double]"0.06052232935307845"),
    @("C5", [double]"0.05405000766046579"),
    @("D5", [double]"0.5636986541578768"),
    @("F5", [double]"0.5352774719538452"),
    @("G5", [double]"0.2934969556344409"),
    @("H5", [double]"0.5878752036870768"),
    @("I5", [double]"0.5062998829755584"),
    @("J5", [double]"0.8734598310264916"),
    @("K5", [double]"0.5076980556173233"),
    @("L5", [double]"0.3233201649915313"),
    @("M5", [double]"0.423299355223121"),
    @("N5", [double]"-6.50592752719859e-16"),
    @("O5", [double]"0.0193538137074701"),
    @("Q5", [double]"-6.50592752719859e-16"),
    @("R5", [double]"-0.2134026769921701"),
    @("S5", [double]"0.4211971398087573"),
    @("T5", [double]"0.2304209406059521"),
    @("U5", [double]"0.2146730676965272"),
    @("B6", [double]"-0.03737587636125757"),
    @("C6", [double]"0.1691006458289873"),
    @("D6", [double]"0.1553023117281974"),
    @("E6", [double]"0.5352774719538452"),
    @("G6", [double]"-0.06193794630210222"),
    @("H6", [double]"0.1665739116006598"),
    @("I6", [double]"0.05360390469665118"),
    @("J6", [double]"0.4975385861924605"),
    @("K6", [double]"0.9245440272241831"),
    @("L6", [double]"-0.01608048406506278"),
    @("M6", [double]"0.04182405402050556"),
    @("N6", [double]"-2.086497594274027e-16"),
    @("O6", [double]"0.05686715317942567"),
    @("Q6", [double]"4.344669041241487e-16"),
    @("R6", [double]"-0.07056329520308711"),
    @("S6", [double]"0.5368144659355543"),
    @("T6", [double]"0.1137821212413632"),
    @("U6", [double]"0.3149986221900413"),
    @("B7", [double]"0.003932712359607337"),
    @("C7", [double]"0.1131570554671354"),
    @("D7", [double]"0.6735708521826163"),
    @("E7", [double]"0.2934969556344409"),
    @("F7", [double]"-0.06193794630210222"),
    @("H7", [double]"0.6790844005770853"),
    @("I7", [double]"0.6778749392792621"),
    @("J7", [double]"0.4261778898454034"),
    @("K7", [double]"0.01219745069401997"),
    @("L7", [double]"0.6530634112319517"),
    @("M7", [double]"0.5393925801821926"),
    @("N7", [double]"-4.581468440064092e-17"),
    @("O7", [double]"-0.1469901637461705"),
    @("Q7", [double]"3.207027908044864e-16"),
    @("R7", [double]"0.08269900451839834"),
    @("S7", [double]"-0.2122788154450698"),
    @("T7", [double]"0.1279721568778777"),
    @("U7", [double]"-0.3038958665051647"),
    @("B8", [double]"0.03251026699987748"),
    @("C8", [double]"0.03464555403083762"),
    @("D8", [double]"0.9869596445666107"),
    @("E8", [double]"0.5878752036870768"),
    @("F8", [double]"0.1665739116006598"),
    @("G8", [double]"0.6790844005770853"),
    @("I8", [double]"0.7720280899487464"),
    @("J8", [double]"0.6802336907754714"),
    @("K8", [double]"0.2746299001602014"),
    @("L8", [double]"0.4493724078540286"),
    @("M8", [double]"0.6074752417321293"),
    @("N8", [double]"1.84520271996319e-16"),
    @("O8", [double]"-0.1514037574916817"),
    @("Q8", [double]"2.039434585222473e-16"),
    @("R8", [double]"-0.1157636184914509"),
    @("S8", [double]"0.06543036595919888"),
    @("T8", [double]"0.2347531325952004"),
    @("U8", [double]"-0.1384112729372551"),
    @("B9", [double]"0.1455275192583088"),
    @("C9", [double]"-0.05464967430273746"),
    @("D9", [double]"0.7703047989162173"),
    @("E9", [double]"0.5062998829755584"),
    @("F9", [double]"0.05360390469665118"),
    @("G9", [double]"0.6778749392792621"),
    @("H9", [double]"0.7720280899487464"),
    @("J9", [double]"0.6373355037709021"),
    @("K9", [double]"0.1313961368501537"),
    @("L9", [double]"0.6854381694108397"),
    @("M9", [double]"0.9453236016231191"),
    @("N9", [double]"4.056173161152782e-15"),
    @("O9", [double]"-0.2037178193153648"),
    @("Q9", [double]"4.069090910073651e-15"),
    @("R9", [double]"-0.1838509992076941"),
    @("S9", [double]"-0.04724653416730241"),
    @("T9", [double]"0.1546249210996476"),
    @("U9", [double]"-0.1572569478461058"),
    @("B10", [double]"0.03110426779743944"),
    @("C10", [double]"-0.02353409310941642"),
    @("D10", [double]"0.6657197040299973"),
    @("E10", [double]"0.8734598310264916"),
    @("F10", [double]"0.4975385861924605"),
    @("G10", [double]"0.4261778898454034"),
    @("H10", [double]"0.6802336907754714"),
    @("I10", [double]"0.6373355037709021"),
    @("K10", [double]"0.5587315479149167"),
    @("L10", [double]"0.3727721876074929"),
    @("M10", [double]"0.5042955498907329"),
    @("N10", [double]"6.897019909926239e-16"),
    @("O10", [double]"-0.0121249135178753"),
    @("Q10", [double]"7.267827431965283e-16"),
    @("R10", [double]"-0.2253718962090097"),
    @("S10", [double]"0.376384401232258"),
    @("T10", [double]"0.2665656676222631"),
    @("U10", [double]"0.1703258938105085"),
    @("B11", [double]"-0.06786464654325448"),
    @("C11", [double]"0.1086358250565714"),
    @("D11", [double]"0.2524737509893859"),
    @("E11", [double]"0.5076980556173233"),
    @("F11", [double]"0.9245440272241831"),
    @("G11", [double]"0.01219745069401997"),
    @("H11", [double]"0.2746299001602014"),
    @("I11", [double]"0.1313961368501537"),
    @("J11", [double]"0.5587315479149167"),
    @("L11", [double]"-0.0506349014461609"),
    @("M11", [double]"0.08617472307857331"),
    @("N11", [double]"-3.503331779890207e-16"),
    @("O11", [double]"0.07611071914631848"),
    @("Q11", [double]"3.108590170888493e-16"),
    @("R11", [double]"-0.06345688522900998"),
    @("S11", [double]"0.4623819073833715"),
    @("T11", [double]"0.1399306648904478"),
    @("U11", [double]"0.2210165422892613"),
    @("B12", [double]"0.09694398396020938"),
    @("C12", [double]"0.1136429820580396"),
    @("D12", [double]"0.4706509334439072"),
    @("E12", [double]"0.3233201649915313"),
    @("F12", [double]"-0.01608048406506278"),
    @("G12", [double]"0.6530634112319517"),
    @("H12", [double]"0.4493724078540286"),
    @("I12", [double]"0.6854381694108397"),
    @("J12", [double]"0.3727721876074929"),
    @("K12", [double]"-0.0506349014461609"),
    @("M12", [double]"0.7095625596926839"),
    @("N12", [double]"3.926509632699553e-17"),
    @("O12", [double]"-0.07155333816120758"),
    @("Q12", [double]"-1.570603853079821e-16"),
    @("R12", [double]"-0.02872318328931278"),
    @("S12", [double]"-0.001769315014088405"),
    @("T12", [double]"0.0491481588542904"),
    @("U12", [double]"-0.08750535152529236"),
    @("B13", [double]"0.1693593265268756"),
    @("C13", [double]"-0.03672462593629056"),
    @("D13", [double]"0.6081882455058323"),
    @("E13", [double]"0.423299355223121"),
    @("F13", [double]"0.04182405402050556"),
    @("G13", [double]"0.5393925801821926"),
    @("H13", [double]"0.6074752417321293"),
    @("I13", [double]"0.9453236016231191"),
    @("J13", [double]"0.5042955498907329"),
    @("K13", [double]"0.08617472307857331"),
    @("L13", [double]"0.7095625596926839"),
    @("N13", [double]"2.198470513904755e-17"),
    @("O13", [double]"-0.2125738933008336"),
    @("Q13", [double]"-5.129764532444429e-17"),
    @("R13", [double]"-0.1875333398232258"),
    @("S13", [double]"-0.05266050773362523"),
    @("T13", [double]"0.0662380697599478"),
    @("U13", [double]"-0.1156739560277123"),
    @("D14", [double]"3.566671208812168e-16"),
    @("E14", [double]"-6.50592752719859e-16"),
    @("F14", [double]"-2.086497594274027e-16"),
    @("G14", [double]"-4.581468440064092e-17"),
    @("H14", [double]"1.84520271996319e-16"),
    @("I14", [double]"4.056173161152782e-15"),
    @("J14", [double]"6.897019909926239e-16"),
    @("K14", [double]"-3.503331779890207e-16"),
    @("L14", [double]"3.926509632699553e-17"),
    @("M14", [double]"2.198470513904755e-17"),
    @("D15", [double]"-0.1375382841619496"),
    @("E15", [double]"0.0193538137074701"),
    @("F15", [double]"0.05686715317942567"),
    @("G15", [double]"-0.1469901637461705"),
    @("H15", [double]"-0.1514037574916817"),
    @("I15", [double]"-0.2037178193153648"),
    @("J15", [double]"-0.0121249135178753"),
    @("K15", [double]"0.07611071914631848"),
    @("L15", [double]"-0.07155333816120758"),
    @("M15", [double]"-0.2125738933008336"),
    @("D17", [double]"4.048653804597595e-16"),
    @("E17", [double]"-6.50592752719859e-16"),
    @("F17", [double]"4.344669041241487e-16"),
    @("G17", [double]"3.207027908044864e-16"),
    @("H17", [double]"2.039434585222473e-16"),
    @("I17", [double]"4.069090910073651e-15"),
    @("J17", [double]"7.267827431965283e-16"),
    @("K17", [double]"3.108590170888493e-16"),
    @("L17", [double]"-1.570603853079821e-16"),
    @("M17", [double]"-5.129764532444429e-17"),
    @("D18", [double]"-0.1123992994262411"),
    @("E18", [double]"-0.2134026769921701"),
    @("F18", [double]"-0.07056329520308711"),
    @("G18", [double]"0.08269900451839834"),
    @("H18", [double]"-0.1157636184914509"),
    @("I18", [double]"-0.1838509992076941"),
    @("J18", [double]"-0.2253718962090097"),
    @("K18", [double]"-0.06345688522900998"),
    @("L18", [double]"-0.02872318328931278"),
    @("M18", [double]"-0.1875333398232258"),
    @("D19", [double]"0.07052404406482074"),
    @("E19", [double]"0.4211971398087573"),
    @("F19", [double]"0.5368144659355543"),
    @("G19", [double]"-0.2122788154450698"),
    @("H19", [double]"0.06543036595919888"),
    @("I19", [double]"-0.04724653416730241"),
    @("J19", [double]"0.376384401232258"),
    @("K19", [double]"0.4623819073833715"),
    @("L19", [double]"-0.001769315014088405"),
    @("M19", [double]"-0.05266050773362523"),
    @("D20", [double]"0.2153517673687143"),
    @("E20", [double]"0.2304209406059521"),
    @("F20", [double]"0.1137821212413632"),
    @("G20", [double]"0.1279721568778777"),
    @("H20", [double]"0.2347531325952004"),
    @("I20", [double]"0.1546249210996476"),
    @("J20", [double]"0.2665656676222631"),
    @("K20", [double]"0.1399306648904478"),
    @("L20", [double]"0.0491481588542904"),
    @("M20", [double]"0.0662380697599478"),
    @("D21", [double]"-0.1390079924904912"),
    @("E21", [double]"0.2146730676965272"),
    @("F21", [double]"0.3149986221900413"),
    @("G21", [double]"-0.3038958665051647"),
    @("H21", [double]"-0.1384112729372551"),
    @("I21", [double]"-0.1572569478461058"),
    @("J21", [double]"0.1703258938105085"),
    @("K21", [double]"0.2210165422892613"),
    @("L21", [double]"-0.08750535152529236"),
    @("M21", [double]"-0.1156739560277123"),
)

foreach ($u in $cellUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}
